$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for each data row (2-33).
# Update every row whose value is 46074 to 46075 (one day later).
for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
